$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.241.05'
$ws.Range('E2').Value = '  -1.59%  '
$ws.Range('D3').Value = '2.251.05'
$ws.Range('E3').Value = '  -1.64%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '247.32'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.76%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.624'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.84%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '77.59'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +5.99%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('E9').Value = '  -4.36%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '41.89'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +7.67%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0954'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -2.21%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '7.12'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -4.19%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.103'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -2.25%  '
$ws.Range('D14').Value = '2.582.33'
$ws.Range('E14').Value = '  -1.85%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '14.81'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -2.78%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.862'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -0.66%  '
$ws.Range('D17').Value = '2.245.30'
$ws.Range('E17').Value = '  -1.75%  '
$ws.Range('D18').Value = '42.098.35'
$ws.Range('E18').Value = '  -1.67%  '
$ws.Range('D19').Value = '0.0₃0986'
$ws.Range('E19').Value = '  -2.07%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '6.12'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -2.53%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '71.96'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.80%  '
$ws.Range('E22').Value = '  +2.15%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '232.32'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -1.91%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '11.28'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -2.27%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '3.61'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -7.69%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.30'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -4.43%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '7.42'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +15.80%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.16'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +1.23%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '169.64'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +1.68%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '20.62'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -2.06%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '33.01'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +6.10%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0831'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +1.45%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.121'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -4.93%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.126'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -0.13%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '4.55'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.74%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '4.94'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +3.06%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0304'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -1.29%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '14.29'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +3.15%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '5.92'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +0.30%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.19'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -6.42%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '113.87'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +10.98%  '
$ws.Range('E43').Value = '  -6.05%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '61.31'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -2.40%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '8.68'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -5.01%  '
$ws.Range('E46').Value = '  -3.53%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.998'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -0.34%  '
$ws.Range('E48').Value = '  -2.42%  '
$ws.Range('B49').Value = 'WOONetwork'
$ws.Range('C49').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.453'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +16.34%  '
$ws.Range('B50').Value = 'TrustWalletToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.17'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -0.65%  '
$ws.Range('E51').Value = '  +0.87%  '
